$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Devices": update a handful of existing device records and
# append a new (test) light device row.
# ---------------------------------------------------------------
$devices = $wb.Worksheets.Item("Devices")

# Row 2 - LI001: rename, clear brand/model, adjust auto on/off values, bump UPDATED_TS
$devices.Range("C2").Value = "Basement_Studio_Test"
$devices.Range("D2").Value = "N/A"
$devices.Range("E2").Value = "N/A"
$devices.Range("G2").Value = 1400.0
$devices.Range("H2").Value = 1400.0
$devices.Range("K2").Value = "2025-07-04T11:14:03.559928237+02:00[Europe/Amsterdam]"

# Row 3 - LI002: just bump UPDATED_TS
$devices.Range("K3").Value = "2025-07-05T00:37:46.157958295+02:00[Europe/Amsterdam]"

# Row 16 - LI015: set brand/model, enable auto-enable, adjust on/off values, bump UPDATED_TS
$devices.Range("D16").Value = "Philips"
$devices.Range("E16").Value = "PLB10025"
$devices.Range("F16").Value = $true
$devices.Range("G16").Value = 1400.0
$devices.Range("H16").Value = 1400.0
$devices.Range("K16").Value = "2025-07-05T00:37:46.377879102+02:00[Europe/Amsterdam]"

# Row 26 - TH001 (thermostat): rename, normalize actions, adjust auto on/off, bump UPDATED_TS
$devices.Range("C26").Value = "Bedroom_Thermostat_Test"
$devices.Range("F26").Value = $false
$devices.Range("G26").Value = 21.0
$devices.Range("H26").Value = 27.0
$devices.Range("I26").Value = "on, off"
$devices.Range("K26").Value = "2025-07-04T11:14:03.794192170+02:00[Europe/Amsterdam]"

# Row 37 - DR001 (dryer): rename, normalize actions, adjust auto on/off, bump UPDATED_TS
$devices.Range("C37").Value = "BasementDryer_Test"
$devices.Range("F37").Value = $false
$devices.Range("G37").Value = 1300.0
$devices.Range("H37").Value = 1300.0
$devices.Range("I37").Value = "on, off"
$devices.Range("K37").Value = "2025-07-04T11:14:04.016417005+02:00[Europe/Amsterdam]"

# Row 38 (new) - a test light device
$devices.Range("A38").Value = "LIGHT"
$devices.Range("B38").Value = "4 errors to this class:LGHT001"
$devices.Range("C38").Value = "MyLight"
$devices.Range("D38").Value = "N/A"
$devices.Range("E38").Value = "N/A"
$devices.Range("F38").Value = "FALSE"
$devices.Range("G38").Value = 400.0
$devices.Range("H38").Value = 400.0
$devices.Range("I38").Value = "on, off"
$devices.Range("J38").Value = "2025-07-04T18:42:46.372028042Z"
$devices.Range("K38").Value = "2025-07-04T18:42:46.372028042Z"
$devices.Range("L38").Value = "N/A"

# ---------------------------------------------------------------
# Sheet "Sensors": update default value + UPDATED_TS on the one row.
# ---------------------------------------------------------------
$sensors = $wb.Worksheets.Item("Sensors")
$sensors.Range("E2").Value = 1200.0
$sensors.Range("H2").Value = "2025-07-04T09:19:26.921304077+02:00[Europe/Amsterdam]"

# ---------------------------------------------------------------
# Sheet "Sense_Control": append a new slave/sensor link row.
# ---------------------------------------------------------------
$senseControl = $wb.Worksheets.Item("Sense_Control")
$senseControl.Range("A3").Value = "LIGHT"
$senseControl.Range("B3").Value = "LI015"
$senseControl.Range("C3").Value = "LIGHT"
$senseControl.Range("D3").Value = "LIs001"
$senseControl.Range("E3").Value = 1400.0
$senseControl.Range("F3").Value = 1400.0
